$d = $word.ActiveDocument

# -------------------------------------------------------------------
# Hunk 1: the "SUN Oct 08" / " 14:29:50 PDT 2017" runs are merged into
# a single run with the combined text.
# -------------------------------------------------------------------
$mergeRange = $d.Content
$mergeRange.Find.Execute("SUN Oct 08 14:29:50 PDT 2017", $true, $false, $false, $false, $false, $true, 1, $false, "SUN Oct 08 14:29:50 PDT 2017", 2)

# -------------------------------------------------------------------
# Hunk 2: a whole new "purchase" block is appended right after the
# "Amount balance ... - 14023.0" paragraph (the unique amount makes it
# easy to locate), and before the following (pre-existing) blank
# paragraph.
# -------------------------------------------------------------------
$anchorRange = $d.Content
$anchorRange.Find.Execute("- 14023.0", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$anchorIndex = $anchorRange.Paragraphs(1).Index

# The paragraph right after the anchor is where we splice in the new
# content -- everything gets inserted *before* it.
$insertionPoint = $d.Paragraphs($anchorIndex + 1).Range

for ($i = 0; $i -lt 11; $i++) {
    $insertionPoint.InsertParagraphBefore()
}

$base = $anchorIndex  # the "Amount balance / - 14023.0" paragraph

# Paragraph base+1: "MON Oct 9" / " 11:05:12 PDT 2017" -- two runs with
# *identical* formatting that must stay split. Word normally coalesces
# adjacent same-format runs typed into the same paragraph, so build the
# two pieces in separate paragraphs and then merge the paragraphs by
# deleting the paragraph mark between them -- that keeps the runs
# distinct.
$p1 = $d.Paragraphs($base + 1).Range
$p1.InsertAfter("MON Oct 9")
$p1After = $d.Paragraphs($base + 1).Range
$p1After.InsertParagraphAfter()
$p2 = $d.Paragraphs($base + 2).Range
$p2.InsertAfter(" 11:05:12 PDT 2017")
$boundary = $d.Paragraphs($base + 1).Range.End - 1
$d.Range($boundary, $boundary + 1).Delete()

# Paragraph base+2: "Person Name" row
$d.Paragraphs($base + 2).Range.InsertAfter("Person Name`t`t`t`t- HSJ")

# Paragraph base+3: divider line
$d.Paragraphs($base + 3).Range.InsertAfter("---------------------------------------------------------------")

# Paragraph base+4: "Item Name" row
$d.Paragraphs($base + 4).Range.InsertAfter("Item Name`t`t`t`t- CARROT")

# Paragraph base+5: "Number of Pockets" row
$d.Paragraphs($base + 5).Range.InsertAfter("Number of Pockets`t`t`t- 1")

# Paragraph base+6: "Number of KGs" row
$d.Paragraphs($base + 6).Range.InsertAfter("Number of KGs`t`t`t- 81")

# Paragraph base+7: "Rate" row
$d.Paragraphs($base + 7).Range.InsertAfter("Rate`t`t`t`t`t- 20")

# Paragraph base+8: "Total Price" row
$d.Paragraphs($base + 8).Range.InsertAfter("Total Price`t`t`t`t- 1620.0")

# Paragraph base+9: "Amount balance" row -- bold
$p9 = $d.Paragraphs($base + 9)
$p9.Range.InsertAfter("Amount balance`t`t`t- 15643.0")
$p9.Range.Font.Bold = $true
$p9.Range.ParagraphFormat.Borders.Application | Out-Null

# Paragraph base+10: blank (not bold)

# Paragraph base+11: blank, bold
$d.Paragraphs($base + 11).Range.Font.Bold = $true

Write-Host "done"
